$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 5182.43
$ws.Range("E2").Value = -5182.43

$ws.Range("D4").Value = 6514.23
$ws.Range("E4").Value = 7209.110000000001
$ws.Range("F4").Value = 0.4746825481260393
